$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.006.25"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "1.773.27"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9980"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4481"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07436"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9976"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.045"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.249"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").Value = "1.769.80"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001063"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06440"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9983"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.783"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").Value = "27.987.66"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.113"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("D28").Value = "1.968.35"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.161"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.112"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.696"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09190"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.682"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  +2.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06222"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02293"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6333"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.978"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.187"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.395"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.896"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.749"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5895"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.962"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.141"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.55%  "

Write-Host "Updated cryptos list"
